$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1437.9
$ws.Range("J19").Value = 1875.8
$ws.Range("L19").Value = 1875.8
$ws.Range("N19").Value = -2225.8

$ws.Range("H39").Value = 87
$ws.Range("I39").Value = 87
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 261
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 35
$ws.Range("N39").ClearContents()

$ws.Range("H40").Value = 2339.4443
$ws.Range("I40").Value = 2243.5
$ws.Range("J40").Value = 2675.25
$ws.Range("K40").Value = 2243.5
$ws.Range("L40").Value = 2675.25
$ws.Range("M40").Value = -2068.5
$ws.Range("N40").Value = -3025.25

$ws.Range("H74").Value = 4267.3335
$ws.Range("I74").Value = 3775
$ws.Range("J74").Value = 5252
$ws.Range("K74").Value = 3775
$ws.Range("L74").Value = 5252
$ws.Range("M74").Value = -2839
$ws.Range("N74").Value = -7124

$ws.Range("H77").Value = 4267.3335
$ws.Range("I77").Value = 3775
$ws.Range("J77").Value = 5252
$ws.Range("K77").Value = 18875
$ws.Range("L77").Value = 26260
$ws.Range("M77").Value = -14195
$ws.Range("N77").Value = -35620

$ws.Range("H113").Value = 94216.82000000001
$ws.Range("I113").Value = 103278.5
$ws.Range("K113").Value = 103278.5
$ws.Range("M113").Value = -100024.5

$ws.Range("H137").Value = 27028378
$ws.Range("I137").Value = 45455596
$ws.Range("J137").Value = 1792.2667
$ws.Range("K137").Value = 136366788
$ws.Range("L137").Value = 5376.800099999999
$ws.Range("M137").Value = -136364238
$ws.Range("N137").Value = -10476.8001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 74119.57000000001
$ws.Range("I2").Value = 113219.336
$ws.Range("J2").Value = 3740
$ws.Range("K2").Value = 113219.336
$ws.Range("L2").Value = 3740
$ws.Range("M2").Value = -113106.336
$ws.Range("N2").Value = -3966

$ws.Range("H32").Value = 23297.176
$ws.Range("I32").Value = 3474.2666
$ws.Range("J32").Value = 171969
$ws.Range("K32").Value = 3474.2666
$ws.Range("L32").Value = 171969
$ws.Range("M32").Value = -3187.2666
$ws.Range("N32").Value = -172543

$ws.Range("H46").Value = 6647.75
$ws.Range("J46").Value = 6647.75
$ws.Range("L46").Value = 6647.75
$ws.Range("N46").Value = -7285.75

$ws.Range("H61").Value = 2212.3096
$ws.Range("I61").Value = 1549.2188
$ws.Range("J61").Value = 4334.2
$ws.Range("K61").Value = 1549.2188
$ws.Range("L61").Value = 4334.2
$ws.Range("M61").Value = -1337.2188
$ws.Range("N61").Value = -4758.2

$ws.Range("H74").Value = 3373.93
$ws.Range("I74").Value = 941.36365
$ws.Range("J74").Value = 11607.23
$ws.Range("K74").Value = 941.36365
$ws.Range("L74").Value = 11607.23
$ws.Range("M74").Value = -67.36365000000001
$ws.Range("N74").Value = -13355.23

$ws.Range("H77").Value = 3373.93
$ws.Range("I77").Value = 941.36365
$ws.Range("J77").Value = 11607.23
$ws.Range("K77").Value = 4706.81825
$ws.Range("L77").Value = 58036.14999999999
$ws.Range("M77").Value = -338.8182500000003
$ws.Range("N77").Value = -66772.14999999999

$ws.Range("H88").Value = 5626.625
$ws.Range("I88").Value = 2549.5
$ws.Range("J88").Value = 6652.3335
$ws.Range("K88").Value = 2549.5
$ws.Range("L88").Value = 6652.3335
$ws.Range("M88").Value = -2143.5
$ws.Range("N88").Value = -7464.3335

$ws.Range("H91").Value = 5626.625
$ws.Range("I91").Value = 2549.5
$ws.Range("J91").Value = 6652.3335
$ws.Range("K91").Value = 2549.5
$ws.Range("L91").Value = 6652.3335
$ws.Range("M91").Value = -1145.5
$ws.Range("N91").Value = -9460.333500000001

$ws.Range("H97").Value = 25649186
$ws.Range("I97").Value = 25649186
$ws.Range("K97").Value = 25649186
$ws.Range("M97").Value = -25648690

$ws.Range("H116").Value = 74119.57000000001
$ws.Range("I116").Value = 113219.336
$ws.Range("J116").Value = 3740
$ws.Range("K116").Value = 113219.336
$ws.Range("L116").Value = 3740
$ws.Range("M116").Value = -110925.336
$ws.Range("N116").Value = -8328

$ws.Range("H122").Value = 1994.5
$ws.Range("I122").Value = 1949.619
$ws.Range("K122").Value = 5848.857
$ws.Range("M122").Value = -3398.857

$ws.Range("H136").Value = 2212.3096
$ws.Range("I136").Value = 1549.2188
$ws.Range("J136").Value = 4334.2
$ws.Range("K136").Value = 4647.6564
$ws.Range("L136").Value = 13002.6
$ws.Range("M136").Value = -2097.6564
$ws.Range("N136").Value = -18102.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 74119.57000000001
$ws.Range("I3").Value = 113219.336
$ws.Range("J3").Value = 3740
$ws.Range("K3").Value = 113219.336
$ws.Range("L3").Value = 3740
$ws.Range("M3").Value = -113105.336
$ws.Range("N3").Value = -3968

$ws.Range("H20").Value = 3333.5386
$ws.Range("I20").Value = 3317.3333
$ws.Range("K20").Value = 3317.3333
$ws.Range("M20").Value = -3070.3333

$ws.Range("H94").Value = 1124.8966
$ws.Range("J94").Value = 1975.5714
$ws.Range("L94").Value = 1975.5714
$ws.Range("N94").Value = -2877.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1703.725
$ws.Range("I31").Value = 1049.6897
$ws.Range("J31").Value = 3428
$ws.Range("K31").Value = 1049.6897
$ws.Range("L31").Value = 3428
$ws.Range("M31").Value = -754.6896999999999
$ws.Range("N31").Value = -4018

$ws.Range("H34").Value = 1703.725
$ws.Range("I34").Value = 1049.6897
$ws.Range("J34").Value = 3428
$ws.Range("K34").Value = 1049.6897
$ws.Range("L34").Value = 3428
$ws.Range("M34").Value = -847.6896999999999
$ws.Range("N34").Value = -3832

$ws.Range("H62").Value = 20216.5
$ws.Range("I62").Value = 31842.715
$ws.Range("J62").Value = 3939.8
$ws.Range("K62").Value = 31842.715
$ws.Range("L62").Value = 3939.8
$ws.Range("M62").Value = -31218.715
$ws.Range("N62").Value = -5187.8

$ws.Range("H65").Value = 20216.5
$ws.Range("I65").Value = 31842.715
$ws.Range("J65").Value = 3939.8
$ws.Range("K65").Value = 159213.575
$ws.Range("L65").Value = 19699
$ws.Range("M65").Value = -156093.575
$ws.Range("N65").Value = -25939

$ws.Range("H122").Value = 2735.1
$ws.Range("I122").Value = 1325
$ws.Range("K122").Value = 3975
$ws.Range("M122").Value = -1525

$ws.Range("H132").Value = 2310.1282
$ws.Range("I132").Value = 1921.7587
$ws.Range("J132").Value = 3436.4
$ws.Range("K132").Value = 5765.2761
$ws.Range("L132").Value = 10309.2
$ws.Range("M132").Value = -3235.2761
$ws.Range("N132").Value = -15369.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 408.93332
$ws.Range("I107").Value = 410.83334
$ws.Range("J107").Value = 401.33334
$ws.Range("K107").Value = 1232.50002
$ws.Range("L107").Value = 1204.00002
$ws.Range("M107").Value = 687.4999800000001
$ws.Range("N107").Value = -5044.000019999999

$ws.Range("H118").Value = 2514.5
$ws.Range("I118").Value = 1029
$ws.Range("J118").Value = 4000
$ws.Range("K118").Value = 3087
$ws.Range("L118").Value = 12000
$ws.Range("M118").Value = -1844
$ws.Range("N118").Value = -14486

$ws.Range("H141").Value = 3076.875
$ws.Range("I141").Value = 3159.2856
$ws.Range("K141").Value = 9477.856800000001
$ws.Range("M141").Value = -4297.856800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I80").Value = 2779
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2779
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1781
$ws.Range("N80").Value = -4996

$ws.Range("I83").Value = 2779
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 13895
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -8903
$ws.Range("N83").Value = -24984

$ws.Range("H97").Value = 903.5294
$ws.Range("I97").Value = 726.9231
$ws.Range("J97").Value = 1477.5
$ws.Range("K97").Value = 726.9231
$ws.Range("L97").Value = 1477.5
$ws.Range("M97").Value = -230.9231
$ws.Range("N97").Value = -2469.5

$ws.Range("H123").Value = 9737.799999999999
$ws.Range("J123").Value = 9737.799999999999
$ws.Range("L123").Value = 9737.799999999999
$ws.Range("N123").Value = -14637.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -888
$ws.Range("N2").ClearContents()

$ws.Range("H68").Value = 2554.7222
$ws.Range("I68").Value = 2301.2
$ws.Range("J68").Value = 2652.2307
$ws.Range("K68").Value = 2301.2
$ws.Range("L68").Value = 2652.2307
$ws.Range("M68").Value = -1552.2
$ws.Range("N68").Value = -4150.2307

$ws.Range("H71").Value = 2554.7222
$ws.Range("I71").Value = 2301.2
$ws.Range("J71").Value = 2652.2307
$ws.Range("K71").Value = 11506
$ws.Range("L71").Value = 13261.1535
$ws.Range("M71").Value = -7762
$ws.Range("N71").Value = -20749.1535

$ws.Range("H97").Value = 27000
$ws.Range("J97").Value = 27000
$ws.Range("L97").Value = 27000
$ws.Range("N97").Value = -28982

$ws.Range("H122").Value = 3095
$ws.Range("I122").Value = 1911.4286
$ws.Range("J122").Value = 3612.8125
$ws.Range("K122").Value = 5734.2858
$ws.Range("L122").Value = 10838.4375
$ws.Range("M122").Value = -3284.2858
$ws.Range("N122").Value = -15738.4375

$ws.Range("H132").Value = 7204.2334
$ws.Range("I132").Value = 8656
$ws.Range("J132").Value = 4696.636
$ws.Range("K132").Value = 25968
$ws.Range("L132").Value = 14089.908
$ws.Range("M132").Value = -23438
$ws.Range("N132").Value = -19149.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 44474444
$ws.Range("I2").Value = 57161428
$ws.Range("J2").Value = 70003
$ws.Range("K2").Value = 57161428
$ws.Range("L2").Value = 70003
$ws.Range("M2").Value = -57161316
$ws.Range("N2").Value = -70227

$ws.Range("H81").Value = 515849.2
$ws.Range("J81").Value = 3779.1667
$ws.Range("L81").Value = 7558.3334
$ws.Range("N81").Value = -9680.3334

$ws.Range("H84").Value = 515849.2
$ws.Range("J84").Value = 3779.1667
$ws.Range("L84").Value = 37791.667
$ws.Range("N84").Value = -48399.667

$ws.Range("H122").Value = 69077.47
$ws.Range("I122").Value = 126900.5
$ws.Range("J122").Value = 2994
$ws.Range("K122").Value = 380701.5
$ws.Range("L122").Value = 8982
$ws.Range("M122").Value = -378251.5
$ws.Range("N122").Value = -13882

$ws.Range("H123").Value = 22483
$ws.Range("J123").Value = 22483
$ws.Range("L123").Value = 22483
$ws.Range("N123").Value = -32283

$ws.Range("H136").Value = 7430770.5
$ws.Range("I136").Value = 7961269.5
$ws.Range("J136").Value = 3788.3333
$ws.Range("K136").Value = 23883808.5
$ws.Range("L136").Value = 11364.9999
$ws.Range("M136").Value = -23881258.5
$ws.Range("N136").Value = -16464.9999
